$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 14495.4
$ws.Range("I10").Value = 10004
$ws.Range("J10").Value = 14994.444
$ws.Range("K10").Value = 10004
$ws.Range("L10").Value = 14994.444
$ws.Range("N10").Value = -15580.444
$ws.Range("M10").Value = -9711

# Row 33
$ws.Range("H33").Value = 628.5
$ws.Range("I33").Value = 259.25
$ws.Range("K33").Value = 259.25
$ws.Range("M33").Value = -30.25

# Row 34
$ws.Range("H34").Value = 11017.833
$ws.Range("I34").Value = 11017.833
$ws.Range("K34").Value = 11017.833
$ws.Range("M34").Value = -10814.833

# Row 36
$ws.Range("H36").Value = 11017.833
$ws.Range("I36").Value = 11017.833
$ws.Range("K36").Value = 11017.833
$ws.Range("M36").Value = -10302.833

# Row 39
$ws.Range("H39").Value = 178.09091
$ws.Range("I39").Value = 78.888885
$ws.Range("J39").Value = 624.5
$ws.Range("K39").Value = 236.666655
$ws.Range("L39").Value = 1873.5
$ws.Range("M39").Value = 59.33334500000001
$ws.Range("N39").Value = -2465.5

# Row 42
$ws.Range("H42").Value = 65.25
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

# Row 132
$ws.Range("H132").Value = 8363.385
$ws.Range("I132").Value = 3373.1
$ws.Range("K132").Value = 10119.3
$ws.Range("M132").Value = -7589.299999999999

# Row 133
$ws.Range("H133").Value = 88748.336
$ws.Range("J133").Value = 88748.336
$ws.Range("L133").Value = 88748.336
$ws.Range("N133").Value = -98868.336

# Row 137
$ws.Range("H137").Value = 2139.8462
$ws.Range("I137").Value = 1332.5
$ws.Range("K137").Value = 3997.5
$ws.Range("M137").Value = -1447.5


$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 3130.875
$ws.Range("I5").Value = 721.2857
$ws.Range("J5").Value = 19998
$ws.Range("K5").Value = 721.2857
$ws.Range("L5").Value = 19998
$ws.Range("M5").Value = -609.2857
$ws.Range("N5").Value = -20222

# Row 32
$ws.Range("H32").Value = 5526.4653
$ws.Range("I32").Value = 2428.3547
$ws.Range("J32").Value = 13529.917
$ws.Range("K32").Value = 2428.3547
$ws.Range("L32").Value = 13529.917
$ws.Range("M32").Value = -2141.3547
$ws.Range("N32").Value = -14103.917

# Row 138
$ws.Range("H138").Value = 99946.336
$ws.Range("J138").Value = 99946.336
$ws.Range("L138").Value = 99946.336
$ws.Range("N138").Value = -110226.336


$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 3130.875
$ws.Range("I4").Value = 721.2857
$ws.Range("J4").Value = 19998
$ws.Range("K4").Value = 721.2857
$ws.Range("L4").Value = 19998
$ws.Range("M4").Value = -606.2857
$ws.Range("N4").Value = -20228

# Row 122
$ws.Range("H122").Value = 99887.5
$ws.Range("J122").Value = 99887.5
$ws.Range("L122").Value = 99887.5
$ws.Range("N122").Value = -109687.5

# Row 134
$ws.Range("H134").Value = 11469.2
$ws.Range("I134").Value = 11604.8
$ws.Range("K134").Value = 34814.39999999999
$ws.Range("M134").Value = -32279.39999999999


$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 3498.5
$ws.Range("I13").Value = 3998
$ws.Range("J13").Value = 2999
$ws.Range("K13").Value = 3998
$ws.Range("L13").Value = 2999
$ws.Range("M13").Value = -3859
$ws.Range("N13").Value = -3277

# Row 64
$ws.Range("H64").Value = 94414.28999999999
$ws.Range("I64").Value = 85450
$ws.Range("K64").Value = 85450
$ws.Range("M64").Value = -85202

# Row 67
$ws.Range("H67").Value = 94414.28999999999
$ws.Range("I67").Value = 85450
$ws.Range("K67").Value = 85450
$ws.Range("M67").Value = -84592

# Row 107
$ws.Range("H107").Value = 590.85297
$ws.Range("I107").Value = 525.5417
$ws.Range("K107").Value = 525.5417
$ws.Range("M107").Value = 1394.4583

# Row 132
$ws.Range("H132").Value = 6046.4614
$ws.Range("I132").Value = 5081.5713
$ws.Range("K132").Value = 15244.7139
$ws.Range("M132").Value = -12714.7139


$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 194.27777
$ws.Range("I11").Value = 88.17646999999999
$ws.Range("J11").Value = 1998
$ws.Range("K11").Value = 264.52941
$ws.Range("L11").Value = 5994
$ws.Range("M11").Value = -124.52941
$ws.Range("N11").Value = -6274

# Row 12
$ws.Range("H12").Value = 98.57143000000001
$ws.Range("I12").Value = 33
$ws.Range("J12").Value = 109.5
$ws.Range("K12").Value = 99
$ws.Range("L12").Value = 328.5
$ws.Range("M12").Value = 74
$ws.Range("N12").Value = -674.5

# Row 121
$ws.Range("H121").Value = 1360.8518
$ws.Range("J121").Value = 1884.7778
$ws.Range("L121").Value = 5654.3334
$ws.Range("N121").Value = -8274.3334


$ws = $wb.Worksheets.Item("GSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()

# Row 14
$ws.Range("H14").Value = 5002747
$ws.Range("I14").Value = 6668330
$ws.Range("J14").Value = 5999
$ws.Range("K14").Value = 6668330
$ws.Range("L14").Value = 5999
$ws.Range("M14").Value = -6668162
$ws.Range("N14").Value = -6335

# Row 116
$ws.Range("H116").Value = 99999
$ws.Range("J116").Value = 99999
$ws.Range("L116").Value = 99999
$ws.Range("N116").Value = -109177

# Row 118
$ws.Range("H118").Value = 99999
$ws.Range("J118").Value = 99999
$ws.Range("L118").Value = 99999
$ws.Range("N118").Value = -103313

# Row 122
$ws.Range("H122").Value = 1270.4286
$ws.Range("I122").Value = 1378.6
$ws.Range("K122").Value = 4135.799999999999
$ws.Range("M122").Value = -1685.799999999999

# Row 132
$ws.Range("H132").Value = 5893.635
$ws.Range("I132").Value = 5192.229
$ws.Range("J132").Value = 8138.1333
$ws.Range("K132").Value = 15576.687
$ws.Range("L132").Value = 24414.3999
$ws.Range("M132").Value = -13046.687
$ws.Range("N132").Value = -29474.3999


$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 10419000
$ws.Range("I46").Value = 20837700
$ws.Range("J46").Value = 299.75
$ws.Range("K46").Value = 20837700
$ws.Range("L46").Value = 299.75
$ws.Range("M46").Value = -20837512
$ws.Range("N46").Value = -675.75


$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 19766.4
$ws.Range("J12").Value = 19766.4
$ws.Range("L12").Value = 19766.4
$ws.Range("N12").Value = -20050.4

# Row 132
$ws.Range("H132").Value = 4361.1816
$ws.Range("I132").Value = 2684.4375
$ws.Range("J132").Value = 8832.5
$ws.Range("K132").Value = 8053.3125
$ws.Range("L132").Value = 26497.5
$ws.Range("M132").Value = -5523.3125
$ws.Range("N132").Value = -31557.5


Write-Output "Applied Zodiark_Profits market data refresh across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets"